$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Row 8 updates: remove the ExisUnits formula, replace with a literal 0,
# enable investment, raise the max investable units, and drop the
# investment cost to reflect the new case study numbers.
$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 3200

# Move the active selection to K11 (bottom-right pane).
[void]$ws.Range("K11").Select()
